$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.033235929787689
$ws.Cells.Item(2, 4).Value = 1.04269230669643
$ws.Cells.Item(2, 5).Value = 1.05163808281508
$ws.Cells.Item(2, 6).Value = 1.056795298537894
$ws.Cells.Item(2, 9).Value = 1.039494626632202
$ws.Cells.Item(2, 10).Value = 1.038361765971013
$ws.Cells.Item(2, 11).Value = 1.045468276411395
$ws.Cells.Item(2, 12).Value = 1.054389030588654
$ws.Cells.Item(2, 13).Value = 1.059532029798143
$ws.Cells.Item(2, 14).Value = 1.016771103568964

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.033963217885662
$ws.Cells.Item(3, 4).Value = 1.043257659258271
$ws.Cells.Item(3, 5).Value = 1.052432034380864
$ws.Cells.Item(3, 6).Value = 1.057534019563753
$ws.Cells.Item(3, 9).Value = 1.039641478269304
$ws.Cells.Item(3, 10).Value = 1.038732810342065
$ws.Cells.Item(3, 11).Value = 1.045845155496561
$ws.Cells.Item(3, 12).Value = 1.054995716485527
$ws.Cells.Item(3, 13).Value = 1.060084651120864
$ws.Cells.Item(3, 14).Value = 1.016894642806203

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.034434533970591
$ws.Cells.Item(4, 4).Value = 1.043624145020133
$ws.Cells.Item(4, 5).Value = 1.052947264268356
$ws.Cells.Item(4, 6).Value = 1.058013187118987
$ws.Cells.Item(4, 9).Value = 1.039735709682618
$ws.Cells.Item(4, 10).Value = 1.038972903136704
$ws.Cells.Item(4, 11).Value = 1.046088980714887
$ws.Cells.Item(4, 12).Value = 1.055389113260234
$ws.Cells.Item(4, 13).Value = 1.060442728390216
$ws.Cells.Item(4, 14).Value = 1.016974563225571

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.034632843700113
$ws.Cells.Item(5, 4).Value = 1.043778372829037
$ws.Cells.Item(5, 5).Value = 1.053164221400292
$ws.Cells.Item(5, 6).Value = 1.058214905979286
$ws.Cells.Item(5, 9).Value = 1.039775134329093
$ws.Cells.Item(5, 10).Value = 1.039073837379775
$ws.Cells.Item(5, 11).Value = 1.046191473685293
$ws.Cells.Item(5, 12).Value = 1.055554694222171
$ws.Cells.Item(5, 13).Value = 1.060593380613003
$ws.Cells.Item(5, 14).Value = 1.017008157060938

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.034666150595437
$ws.Cells.Item(6, 4).Value = 1.043804277517348
$ws.Cells.Item(6, 5).Value = 1.053200670179137
$ws.Cells.Item(6, 6).Value = 1.058248791631227
$ws.Cells.Item(6, 9).Value = 1.039781742722513
$ws.Cells.Item(6, 10).Value = 1.039090784600116
$ws.Cells.Item(6, 11).Value = 1.046208682011498
$ws.Cells.Item(6, 12).Value = 1.055582507490187
$ws.Cells.Item(6, 13).Value = 1.060618682565563
$ws.Cells.Item(6, 14).Value = 1.017013797324056

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.034437183135551
$ws.Cells.Item(7, 4).Value = 1.043626205204481
$ws.Cells.Item(7, 5).Value = 1.052950161870115
$ws.Cells.Item(7, 6).Value = 1.058015881411376
$ws.Cells.Item(7, 9).Value = 1.03973623722515
$ws.Cells.Item(7, 10).Value = 1.038974251829665
$ws.Cells.Item(7, 11).Value = 1.046090350276544
$ws.Cells.Item(7, 12).Value = 1.055391324989567
$ws.Cells.Item(7, 13).Value = 1.060444740956939
$ws.Cells.Item(7, 14).Value = 1.01697501212713

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.033481571555759
$ws.Cells.Item(8, 4).Value = 1.042883231395439
$ws.Cells.Item(8, 5).Value = 1.051906092566975
$ws.Cells.Item(8, 6).Value = 1.057044710064492
$ws.Cells.Item(8, 9).Value = 1.039544419233081
$ws.Cells.Item(8, 10).Value = 1.038487160893309
$ws.Cells.Item(8, 11).Value = 1.045595652061739
$ws.Cells.Item(8, 12).Value = 1.05459388987013
$ws.Cells.Item(8, 13).Value = 1.059718687325565
$ws.Cells.Item(8, 14).Value = 1.016812857623974

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.031803204751315
$ws.Cells.Item(9, 4).Value = 1.041579195740075
$ws.Cells.Item(9, 5).Value = 1.050077821461379
$ws.Cells.Item(9, 6).Value = 1.055342408695926
$ws.Cells.Item(9, 9).Value = 1.039200384398546
$ws.Cells.Item(9, 10).Value = 1.037628919926467
$ws.Cells.Item(9, 11).Value = 1.044723680236104
$ws.Cells.Item(9, 12).Value = 1.053195149252653
$ws.Cells.Item(9, 13).Value = 1.058443153737695
$ws.Cells.Item(9, 14).Value = 1.016527005448728

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.030688141287255
$ws.Cells.Item(10, 4).Value = 1.040713438854014
$ws.Cells.Item(10, 5).Value = 1.048866849520593
$ws.Cells.Item(10, 6).Value = 1.054213744385397
$ws.Cells.Item(10, 9).Value = 1.038967022047901
$ws.Cells.Item(10, 10).Value = 1.037056888658564
$ws.Cells.Item(10, 11).Value = 1.044142280965547
$ws.Cells.Item(10, 12).Value = 1.052267097588147
$ws.Cells.Item(10, 13).Value = 1.057595505138311
$ws.Cells.Item(10, 14).Value = 1.016336387681764

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.030206244801616
$ws.Cells.Item(11, 4).Value = 1.040339435580239
$ws.Cells.Item(11, 5).Value = 1.048344382001881
$ws.Cells.Item(11, 6).Value = 1.053726520157134
$ws.Cells.Item(11, 9).Value = 1.038865033306206
$ws.Cells.Item(11, 10).Value = 1.036809240176841
$ws.Cells.Item(11, 11).Value = 1.043890526322737
$ws.Cells.Item(11, 12).Value = 1.051866318144268
$ws.Cells.Item(11, 13).Value = 1.057229128994575
$ws.Cells.Item(11, 14).Value = 1.016253842211098

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.03002738930533
$ws.Cells.Item(12, 4).Value = 1.040200647648386
$ws.Cells.Item(12, 5).Value = 1.048150600809922
$ws.Cells.Item(12, 6).Value = 1.053545770309299
$ws.Cells.Item(12, 9).Value = 1.038827009515918
$ws.Cells.Item(12, 10).Value = 1.036717260595125
$ws.Cells.Item(12, 11).Value = 1.04379701414246
$ws.Cells.Item(12, 12).Value = 1.051717613955572
$ws.Cells.Item(12, 13).Value = 1.057093141807602
$ws.Cells.Item(12, 14).Value = 1.0162231806179

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.030065747942825
$ws.Cells.Item(13, 4).Value = 1.040230412055187
$ws.Cells.Item(13, 5).Value = 1.048192154530514
$ws.Cells.Item(13, 6).Value = 1.053584531458931
$ws.Cells.Item(13, 9).Value = 1.038835172099148
$ws.Cells.Item(13, 10).Value = 1.036736990139753
$ws.Cells.Item(13, 11).Value = 1.043817072767579
$ws.Cells.Item(13, 12).Value = 1.051749504095408
$ws.Cells.Item(13, 13).Value = 1.057122306910082
$ws.Cells.Item(13, 14).Value = 1.016229757650019

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.030191457635234
$ws.Cells.Item(14, 4).Value = 1.040327960590544
$ws.Cells.Item(14, 5).Value = 1.048328358129817
$ws.Cells.Item(14, 6).Value = 1.053711574675636
$ws.Cells.Item(14, 9).Value = 1.038861893112089
$ws.Cells.Item(14, 10).Value = 1.036801636943037
$ws.Cells.Item(14, 11).Value = 1.043882796554056
$ws.Cells.Item(14, 12).Value = 1.051854022870466
$ws.Cells.Item(14, 13).Value = 1.057217886175476
$ws.Cells.Item(14, 14).Value = 1.016251307722282

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.030268930446278
$ws.Cells.Item(15, 4).Value = 1.040388081218867
$ws.Cells.Item(15, 5).Value = 1.048412315698113
$ws.Cells.Item(15, 6).Value = 1.053789880333926
$ws.Cells.Item(15, 9).Value = 1.038878338198793
$ws.Cells.Item(15, 10).Value = 1.036841469084805
$ws.Cells.Item(15, 11).Value = 1.04392329128441
$ws.Cells.Item(15, 12).Value = 1.051918442012789
$ws.Cells.Item(15, 13).Value = 1.057276789184806
$ws.Cells.Item(15, 14).Value = 1.016264585377137

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.030720143012913
$ws.Cells.Item(16, 4).Value = 1.040738278833649
$ws.Cells.Item(16, 5).Value = 1.048901564018387
$ws.Cells.Item(16, 6).Value = 1.054246111547493
$ws.Cells.Item(16, 9).Value = 1.038973770939984
$ws.Cells.Item(16, 10).Value = 1.03707332533297
$ws.Cells.Item(16, 11).Value = 1.044158989111553
$ws.Cells.Item(16, 12).Value = 1.052293718772437
$ws.Cells.Item(16, 13).Value = 1.057619834409021
$ws.Cells.Item(16, 14).Value = 1.016341865852166

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.031003428094779
$ws.Cells.Item(17, 4).Value = 1.04095818430335
$ws.Cells.Item(17, 5).Value = 1.049208964589614
$ws.Cells.Item(17, 6).Value = 1.054532695203348
$ws.Cells.Item(17, 9).Value = 1.039033381920508
$ws.Cells.Item(17, 10).Value = 1.037218775591881
$ws.Cells.Item(17, 11).Value = 1.044306835832711
$ws.Cells.Item(17, 12).Value = 1.052529408616603
$ws.Cells.Item(17, 13).Value = 1.057835195985354
$ws.Cells.Item(17, 14).Value = 1.01639034039209

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.031168753405084
$ws.Cells.Item(18, 4).Value = 1.041086535747256
$ws.Cells.Item(18, 5).Value = 1.04938844846924
$ws.Cells.Item(18, 6).Value = 1.054699998677498
$ws.Cells.Item(18, 9).Value = 1.039068061121703
$ws.Cells.Item(18, 10).Value = 1.03730361850189
$ws.Cells.Item(18, 11).Value = 1.044393071752197
$ws.Cells.Item(18, 12).Value = 1.052666985919819
$ws.Cells.Item(18, 13).Value = 1.057960876428841
$ws.Cells.Item(18, 14).Value = 1.016418614090632

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.031225140251073
$ws.Cells.Item(19, 4).Value = 1.041130314509649
$ws.Cells.Item(19, 5).Value = 1.049449678688025
$ws.Cells.Item(19, 6).Value = 1.054757069205363
$ws.Cells.Item(19, 9).Value = 1.039079870395101
$ws.Cells.Item(19, 10).Value = 1.03733254841871
$ws.Cells.Item(19, 11).Value = 1.044422475831336
$ws.Cells.Item(19, 12).Value = 1.052713913705458
$ws.Cells.Item(19, 13).Value = 1.058003740956338
$ws.Cells.Item(19, 14).Value = 1.01642825456116

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.030973024971107
$ws.Cells.Item(20, 4).Value = 1.040934581792871
$ws.Cells.Item(20, 5).Value = 1.049175964541581
$ws.Cells.Item(20, 6).Value = 1.054501932569941
$ws.Cells.Item(20, 9).Value = 1.039026995620413
$ws.Cells.Item(20, 10).Value = 1.037203169707832
$ws.Cells.Item(20, 11).Value = 1.044290973327461
$ws.Cells.Item(20, 12).Value = 1.052504110617348
$ws.Cells.Item(20, 13).Value = 1.057812083126186
$ws.Cells.Item(20, 14).Value = 1.016385139599642

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.030154435318683
$ws.Cells.Item(21, 4).Value = 1.040299231275496
$ws.Cells.Item(21, 5).Value = 1.04828824163668
$ws.Cells.Item(21, 6).Value = 1.053674157331808
$ws.Cells.Item(21, 9).Value = 1.038854028315526
$ws.Cells.Item(21, 10).Value = 1.03678259983652
$ws.Cells.Item(21, 11).Value = 1.043863442506062
$ws.Cells.Item(21, 12).Value = 1.051823240177297
$ws.Cells.Item(21, 13).Value = 1.057189737654703
$ws.Cells.Item(21, 14).Value = 1.016244961773907

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.029640579939189
$ws.Cells.Item(22, 4).Value = 1.039900534801826
$ws.Cells.Item(22, 5).Value = 1.047731753821288
$ws.Cells.Item(22, 6).Value = 1.053155016253701
$ws.Cells.Item(22, 9).Value = 1.038744463652549
$ws.Cells.Item(22, 10).Value = 1.036518218347854
$ws.Cells.Item(22, 11).Value = 1.043594641693069
$ws.Cells.Item(22, 12).Value = 1.051396094747153
$ws.Cells.Item(22, 13).Value = 1.056799030943012
$ws.Cells.Item(22, 14).Value = 1.016156823581315

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.029912905554641
$ws.Cells.Item(23, 4).Value = 1.040111817332465
$ws.Cells.Item(23, 5).Value = 1.048026600518817
$ws.Cells.Item(23, 6).Value = 1.053430097423406
$ws.Cells.Item(23, 9).Value = 1.038802622776425
$ws.Cells.Item(23, 10).Value = 1.036658367024233
$ws.Cells.Item(23, 11).Value = 1.043737137155093
$ws.Cells.Item(23, 12).Value = 1.051622442418444
$ws.Cells.Item(23, 13).Value = 1.057006095736349
$ws.Cells.Item(23, 14).Value = 1.016203547413156

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.030986762556471
$ws.Cells.Item(24, 4).Value = 1.040945246491198
$ws.Cells.Item(24, 5).Value = 1.049190875279784
$ws.Cells.Item(24, 6).Value = 1.054515832434969
$ws.Cells.Item(24, 9).Value = 1.039029881595428
$ws.Cells.Item(24, 10).Value = 1.037210221322188
$ws.Cells.Item(24, 11).Value = 1.044298140912755
$ws.Cells.Item(24, 12).Value = 1.052515541375036
$ws.Cells.Item(24, 13).Value = 1.057822526635862
$ws.Cells.Item(24, 14).Value = 1.016387489616052

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.032236433472838
$ws.Cells.Item(25, 4).Value = 1.041915694281897
$ws.Cells.Item(25, 5).Value = 1.05054909515263
$ws.Cells.Item(25, 6).Value = 1.055781410858242
$ws.Cells.Item(25, 9).Value = 1.03929003510362
$ws.Cells.Item(25, 10).Value = 1.037850778737814
$ws.Cells.Item(25, 11).Value = 1.044949126904552
$ws.Cells.Item(25, 12).Value = 1.053555982396828
$ws.Cells.Item(25, 13).Value = 1.058772440502859
$ws.Cells.Item(25, 14).Value = 1.016600915873663
